$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")
$ws.Activate()

# --- Row 92: Multipliers ---
$ws.Cells.Item(92, 1).Value = "Multipliers"
$ws.Cells.Item(92, 2).Value = "Liz Wiseman"
$ws.Cells.Item(92, 3).Value2 = 44000
$ws.Cells.Item(92, 3).NumberFormat = "m/d/yy"
$ws.Cells.Item(92, 4).Value2 = 44003
$ws.Cells.Item(92, 4).NumberFormat = "m/d/yy"
$ws.Cells.Item(92, 5).Value = "leadership;growth mindset;development"
$ws.Cells.Item(92, 6).Value = "Audio"
$ws.Cells.Item(92, 7).Value = "7 Hours 53 Mins"

# --- Row 93: An Economist Walks Into a Brothel ---
$ws.Cells.Item(93, 1).Value = "An Economist Walks Into a Brothel"
$ws.Cells.Item(93, 2).Value = "Allison Schrager"
$ws.Cells.Item(93, 3).Value2 = 44003
$ws.Cells.Item(93, 3).NumberFormat = "m/d/yy"
$ws.Cells.Item(93, 4).Value2 = 44004
$ws.Cells.Item(93, 4).NumberFormat = "m/d/yy"
$ws.Cells.Item(93, 5).Value = "economics;risk;risk management;"
$ws.Cells.Item(93, 6).Value = "Audio"
$ws.Cells.Item(93, 7).Value = "7 Hours 38 Mins"

$ws.Range("A94").Select()
$excel.ActiveWindow.ScrollRow = 69
$excel.ActiveWindow.ScrollColumn = 1
